# Atualizado por script em 02-12-2023 14:46
# Append the newest match row (Tuzla City 3 x 0 Zeljeznicar) to the
# betexplorer odds sheet, matching the style of the preceding row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 91

# Copy the formatting (font, borders, number formats, alignment) of the
# previous data row down onto the new row before filling in values, so the
# new row's style exactly mirrors row 90 (bold/bordered index column,
# datetime-formatted match date column, etc.).
$ws.Range("A90:V90").Copy()
$ws.Range("A91:V91").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 90
$ws.Cells.Item($row, 2).Value = "bosnia-and-herzegovina"
$ws.Cells.Item($row, 3).Value = "premijer-liga-bih"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45262.5625
$ws.Cells.Item($row, 6).Value = "Tuzla City"
$ws.Cells.Item($row, 7).Value = 3
$ws.Cells.Item($row, 8).Value = "Zeljeznicar"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 2.3
$ws.Cells.Item($row, 11).Value = "01/12/2023 01:42"
$ws.Cells.Item($row, 12).Value = 2.16
$ws.Cells.Item($row, 13).Value = "02/12/2023 13:21"
$ws.Cells.Item($row, 14).Value = 3.08
$ws.Cells.Item($row, 15).Value = "01/12/2023 01:42"
$ws.Cells.Item($row, 16).Value = 3.19
$ws.Cells.Item($row, 17).Value = "02/12/2023 13:20"
$ws.Cells.Item($row, 18).Value = 2.89
$ws.Cells.Item($row, 19).Value = "01/12/2023 01:42"
$ws.Cells.Item($row, 20).Value = 3.39
$ws.Cells.Item($row, 21).Value = "02/12/2023 13:24"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/tuzla-city-zeljeznicar/beKsVUMG/"

$wb.Save()
